# Updated cryptos list values (Price / Volume(1h)) per source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text-valued cells (e.g. "29.058.40", "0.9967") are written as text,
# not auto-converted to numbers, by forcing Text number format first.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
}

Set-TextValue $ws.Range("D2") '29.058.40'
Set-TextValue $ws.Range("E2") '  +0.63%  '

Set-TextValue $ws.Range("D3") '1.832.44'
Set-TextValue $ws.Range("E3") '  +0.67%  '

Set-TextValue $ws.Range("D4") '0.9967'
Set-TextValue $ws.Range("E4") '  +0.33%  '

Set-TextValue $ws.Range("D5") '241.64'
Set-TextValue $ws.Range("E5") '  -0.39%  '

Set-TextValue $ws.Range("D6") '0.6183'
Set-TextValue $ws.Range("E6") '  -1.50%  '

Set-TextValue $ws.Range("D7") '0.9997'
Set-TextValue $ws.Range("E7") '  +0.59%  '

Set-TextValue $ws.Range("D8") '0.07450'
Set-TextValue $ws.Range("E8") '  +0.12%  '

Set-TextValue $ws.Range("D9") '0.2932'
Set-TextValue $ws.Range("E9") '  +0.24%  '

Set-TextValue $ws.Range("D10") '23.05'
Set-TextValue $ws.Range("E10") '  +0.42%  '

Set-TextValue $ws.Range("D11") '0.07668'
Set-TextValue $ws.Range("E11") '  +0.10%  '

Set-TextValue $ws.Range("D12") '1.833.57'
Set-TextValue $ws.Range("E12") '  +0.55%  '

Set-TextValue $ws.Range("D13") '4.999'
Set-TextValue $ws.Range("E13") '  +0.61%  '

Set-TextValue $ws.Range("D14") '0.6733'
Set-TextValue $ws.Range("E14") '  +1.36%  '

Set-TextValue $ws.Range("D15") '82.90'
Set-TextValue $ws.Range("E15") '  +0.22%  '

Set-TextValue $ws.Range("D16") '0.000009165'
Set-TextValue $ws.Range("E16") '  -4.73%  '

Set-TextValue $ws.Range("D17") '5.913'
Set-TextValue $ws.Range("E17") '  -1.83%  '

Set-TextValue $ws.Range("D18") '29.046.42'
Set-TextValue $ws.Range("E18") '  +0.53%  '

Set-TextValue $ws.Range("D19") '2.070.47'
Set-TextValue $ws.Range("E19") '  +0.18%  '

Set-TextValue $ws.Range("D20") '239.17'
Set-TextValue $ws.Range("E20") '  +6.34%  '

Set-TextValue $ws.Range("D21") '12.68'
Set-TextValue $ws.Range("E21") '  +1.40%  '

Set-TextValue $ws.Range("D22") '0.9991'
Set-TextValue $ws.Range("E22") '  +0.54%  '

Set-TextValue $ws.Range("D23") '7.205'
Set-TextValue $ws.Range("E23") '  +1.48%  '

Set-TextValue $ws.Range("D24") '0.9981'
Set-TextValue $ws.Range("E24") '  +0.39%  '

Set-TextValue $ws.Range("D25") '158.82'
Set-TextValue $ws.Range("E25") '  -0.60%  '

Set-TextValue $ws.Range("D26") '0.1409'
Set-TextValue $ws.Range("E26") '  +0.31%  '

Set-TextValue $ws.Range("E27") '  +0.48%  '

Set-TextValue $ws.Range("E28") '  +0.15%  '

Set-TextValue $ws.Range("D29") '1.495'
Set-TextValue $ws.Range("E29") '  +0.11%  '

Set-TextValue $ws.Range("D30") '0.05615'
Set-TextValue $ws.Range("E30") '  +3.43%  '

Set-TextValue $ws.Range("B31") 'InternetComputer(DFINITY)'
Set-TextValue $ws.Range("C31") 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue $ws.Range("D31") '4.115'
Set-TextValue $ws.Range("E31") '  +1.87%  '

Set-TextValue $ws.Range("B32") 'Filecoin'
Set-TextValue $ws.Range("C32") 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue $ws.Range("D32") '4.129'
Set-TextValue $ws.Range("E32") '  +0.52%  '

Set-TextValue $ws.Range("D33") '1.202'
Set-TextValue $ws.Range("E33") '  +0.49%  '

Set-TextValue $ws.Range("D34") '1.841'
Set-TextValue $ws.Range("E34") '  -0.31%  '

Set-TextValue $ws.Range("D35") '0.7402'
Set-TextValue $ws.Range("E35") '  +0.22%  '

Set-TextValue $ws.Range("E36") '  +1.08%  '

Set-TextValue $ws.Range("D37") '2.655'
Set-TextValue $ws.Range("E37") '  +1.77%  '

Set-TextValue $ws.Range("E38") '  +1.34%  '

Set-TextValue $ws.Range("D39") '0.01785'
Set-TextValue $ws.Range("E39") '  +0.94%  '

Set-TextValue $ws.Range("D40") '1.212.49'
Set-TextValue $ws.Range("E40") '  -2.30%  '

Set-TextValue $ws.Range("D41") '6.397'
Set-TextValue $ws.Range("E41") '  -3.39%  '

Set-TextValue $ws.Range("D42") '0.8956'
Set-TextValue $ws.Range("E42") '  -0.17%  '

Set-TextValue $ws.Range("D43") '0.9986'
Set-TextValue $ws.Range("E43") '  +0.45%  '

Set-TextValue $ws.Range("D44") '101.45'
Set-TextValue $ws.Range("E44") '  +0.45%  '

Set-TextValue $ws.Range("D45") '1.972.49'
Set-TextValue $ws.Range("E45") '  +0.15%  '

Set-TextValue $ws.Range("D46") '65.40'
Set-TextValue $ws.Range("E46") '  +1.02%  '

Set-TextValue $ws.Range("D47") '0.5077'
Set-TextValue $ws.Range("E47") '  +0.22%  '

Set-TextValue $ws.Range("D48") '0.00000000118'
Set-TextValue $ws.Range("E48") '  -3.27%  '

Set-TextValue $ws.Range("D49") '0.4063'
Set-TextValue $ws.Range("E49") '  +0.96%  '

Set-TextValue $ws.Range("D50") '9.157'
Set-TextValue $ws.Range("E50") '  +2.68%  '

Set-TextValue $ws.Range("D51") '0.05808'
Set-TextValue $ws.Range("E51") '  +0.53%  '

